# Commit 15: JS6 - set() to optional chaining
#
# Inserts two new paragraphs right after the paragraph that ends the
# "Commit 14" bullet ("Sets, has(), for of loop in sets, finding length
# of set, map objects, map iterator"):
#   1) "Commit 15:"
#   2) "set(), get(), keys(), adding key value pair using arrays,
#       extending key value pair in an object using map, Object.assign(),
#       optional chaining"
#
# Both new paragraphs reuse the same paragraph/run formatting
# (ListParagraph style, justified, 12pt/sz 24) as the paragraph they
# follow, which InsertParagraphAfter() naturally inherits.

$d = $word.ActiveDocument

$anchorText = "Sets, has(), for of loop in sets, finding length of set, map objects, map iterator"

$searchRng = $d.Content.Duplicate
$found = $searchRng.Find.Execute($anchorText, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if (-not $found) {
    throw "edit.ps1: could not find anchor paragraph text: $anchorText"
}

$anchorPara = $searchRng.Paragraphs(1)

# Insert the "Commit 15:" paragraph right after the anchor paragraph.
$anchorPara.Range.InsertParagraphAfter()
$commitPara = $anchorPara.Next()
$commitPara.Range.Text = "Commit 15:"

# Insert the description paragraph right after the "Commit 15:" paragraph.
$commitPara.Range.InsertParagraphAfter()
$descPara = $commitPara.Next()
$descPara.Range.Text = "set(), get(), keys(), adding key value pair using arrays, extending key value pair in an object using map, Object.assign(), optional chaining"
